$d = $word.ActiveDocument

# --- Paragraph 13 replacement ---
$xml1 = '<w:p w14:paraId="754F5847" w14:textId="585D4241" w:rsidR="00D6041D" w:rsidRDefault="003A0BD4" w:rsidP="002D6486"><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>project3.exs</w:t></w:r><w:r><w:t xml:space="preserve"> reads the input argument</w:t></w:r><w:r><w:t xml:space="preserve">s and starts the </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>GLOBALSUP</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">We use a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>high level</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> supervisor so that we can terminate once we hear back how many hops all requests took to resolve. </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>GLOBALSUP</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>initializes the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>the</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>MAINPROJ</w:t></w:r><w:r><w:t xml:space="preserve"> module</w:t></w:r><w:r><w:t xml:space="preserve"> with the input arguments</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r w:rsidR="002D6486"><w:t xml:space="preserve"> The </w:t></w:r><w:r w:rsidR="002D6486" w:rsidRPr="002D6486"><w:rPr><w:u w:val="single"/></w:rPr><w:t>MAINPROJ</w:t></w:r><w:r w:rsidR="002D6486"><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="002D6486"><w:t xml:space="preserve">starts the </w:t></w:r><w:r w:rsidR="00D6041D"><w:t xml:space="preserve">dynamic </w:t></w:r><w:r w:rsidR="002D6486"><w:t xml:space="preserve">supervisor </w:t></w:r><w:r w:rsidR="002D6486"><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">TAPESTRY </w:t></w:r><w:r w:rsidR="002D6486"><w:t xml:space="preserve">module, starts all the nodes from </w:t></w:r><w:r w:rsidR="002D6486"><w:rPr><w:u w:val="single"/></w:rPr><w:t>TAPNODE</w:t></w:r><w:r w:rsidR="002D6486"><w:t xml:space="preserve"> module and adds them to the Tapestry mesh. </w:t></w:r><w:r w:rsidR="00ED318E"><w:t xml:space="preserve">Once they are all inserted to the mesh the </w:t></w:r><w:r w:rsidR="00ED318E" w:rsidRPr="002D6486"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t>MAINPROJ</w:t></w:r><w:r w:rsidR="00ED318E"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00ED318E"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">tells each node </w:t></w:r><w:r w:rsidR="00ED318E"><w:rPr><w:rFonts w:eastAsiaTheme="minorHAnsi"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:rPr><w:t>[number of requests]</w:t></w:r><w:r w:rsidR="00ED318E"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> random objects from the Tapestry mesh. </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">As requests are </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t>fulfilled</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> they send the number of hops it took back to </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t>MAINPROJ</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">which adds it to its known number of hops. </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t>Once it receives the number of hops it is expecting (</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorHAnsi"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:rPr><w:t>number of nodes</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorHAnsi"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:rPr><w:t xml:space="preserve"> * </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorHAnsi"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:rPr><w:t>number of requests</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t>) it gets the max</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">outputs it to the terminal and terminates the program. </w:t></w:r></w:p>'
$p1 = $d.Paragraphs(13)
$rng1 = $p1.Range
$rng1.Collapse(1)
$rng1.InsertXML($xml1)

# --- Paragraph 23 replacement ---
$xml2 = '<w:p w14:paraId="767BA22F" w14:textId="50A848CE" w:rsidR="00D85430" w:rsidRDefault="00C74AA9" w:rsidP="002D6486"><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">For N to populate its neighbor map it uses routing. </w:t></w:r><w:r w:rsidR="00E159DF"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">If there is an element where </w:t></w:r><w:r w:rsidR="00BA231C"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t>it would be in B’s routing table it routes to that element</w:t></w:r><w:r w:rsidR="00277ADD"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t>, copies that level that’s its routed too</w:t></w:r><w:r w:rsidR="00BA231C"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> and tries to find a possible closer neighbor. It continues this until no neighbors are available. </w:t></w:r><w:r w:rsidR="005E6282"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">An interesting observation made was </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">that as the network is small it’s more likely you will not have anything in common with the gateway node and that the gateway node will not have anything in common with its neighbors. This leads to very large first levels and an almost fully connected network. </w:t></w:r><w:r w:rsidR="0040753B"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">This is not the case with larger networks as it is more likely that you “match” prefixes with other elements and can better place yourself. </w:t></w:r><w:r w:rsidR="005E6282"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'
$p2 = $d.Paragraphs(23)
$rng2 = $p2.Range
$rng2.Collapse(1)
$rng2.InsertXML($xml2)

# --- Paragraph 25 replacement ---
$xml3 = '<w:p w14:paraId="47E86C7E" w14:textId="63DFACE2" w:rsidR="001C5638" w:rsidRDefault="001C5638" w:rsidP="001C5638"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>ROUTING</w:t></w:r></w:p>
'
$p3 = $d.Paragraphs(25)
$rng3 = $p3.Range
$rng3.Collapse(1)
$rng3.InsertXML($xml3)

# --- Paragraph 26 replacement ---
$xml4 = '<w:p w14:paraId="16F83598" w14:textId="0AF8B5C7" w:rsidR="00AC0BEF" w:rsidRPr="006D6A40" w:rsidRDefault="00B221F3" w:rsidP="00AC0BEF"><w:r><w:t xml:space="preserve">To Route an object </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>MAINPROJ</w:t></w:r><w:r><w:t xml:space="preserve"> tells a node which </w:t></w:r><w:r w:rsidR="00EB0337"><w:t xml:space="preserve">target </w:t></w:r><w:r><w:t>id to route towards</w:t></w:r><w:r w:rsidR="0099584E"><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0099584E"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>R</w:t></w:r><w:r w:rsidR="00A324F8"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>outeToObject</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0099584E"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="004121AB"><w:t xml:space="preserve">finds the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004121AB"><w:t>repfix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004121AB"><w:t xml:space="preserve"> match length to check the neighbor map to see if </w:t></w:r><w:r w:rsidR="006971D8"><w:t>that level</w:t></w:r><w:r w:rsidR="004121AB"><w:t xml:space="preserve"> exists. </w:t></w:r><w:r w:rsidR="006971D8"><w:t xml:space="preserve">If </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="006971D8"><w:t>so</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="006971D8"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00176A23"><w:t xml:space="preserve">it gets that level </w:t></w:r><w:r w:rsidR="00BC7031"><w:t xml:space="preserve">and checks if there is a matching node. </w:t></w:r><w:r w:rsidR="00313F73"><w:t xml:space="preserve">If so that node is a neighbor and you can send a direct message. </w:t></w:r><w:r w:rsidR="006E59C2"><w:t xml:space="preserve">If not the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006E59C2"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>findNetHop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006E59C2"><w:t xml:space="preserve"> is used to find the next neighbor closest to the target id. </w:t></w:r><w:r w:rsidR="006D6A40"><w:t xml:space="preserve">Once </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006D6A40"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>findNextHop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006D6A40"><w:t xml:space="preserve"> finds that neighbor it contacts </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006D6A40"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>nextHop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006D6A40"><w:t xml:space="preserve"> to move to that neighbor and repeat the process checking there. </w:t></w:r></w:p>
'
$p4 = $d.Paragraphs(26)
$rng4 = $p4.Range
$rng4.Collapse(1)
$rng4.InsertXML($xml4)
